# Balancing Fetch Boss & Weapon & Enemy
$wb = $excel.ActiveWorkbook

# --- WeaponDB sheet: rebalance damage/reloadSpeed/magazine/useTime/bulletSpeed ---
$wsWeapon = $wb.Worksheets.Item("WeaponDB")

$weaponRows = @(
    @(2,25,0.7,80,0.07,0),
    @(3,20,1,50,0.12,0),
    @(4,15,1.5,30,0.2,0),
    @(5,10,2,20,0.3,0),
    @(6,40,0.5,20,0.3,0),
    @(7,30,0.7,15,0.5,0),
    @(8,20,1,10,0.8,0),
    @(9,15,1.5,6,1,0),
    @(10,100,1.5,10,1,0),
    @(11,70,2,10,1.3,0),
    @(12,50,2.5,6,1.5,0),
    @(13,30,3,3,1.7,0),
    @(14,21,0.8,2,0.3,0),
    @(15,15,0.8,2,0.3,0),
    @(16,12,1.5,12,0.8,0),
    @(17,7,2,6,1,0)
)

foreach ($row in $weaponRows) {
    $r = $row[0]
    $wsWeapon.Cells.Item($r, 2).Value = $row[1]
    $wsWeapon.Cells.Item($r, 3).Value = $row[2]
    $wsWeapon.Cells.Item($r, 4).Value = $row[3]
    $wsWeapon.Cells.Item($r, 5).Value = $row[4]
    $wsWeapon.Cells.Item($r, 6).Value = $row[5]
}

$wsWeapon.Activate()
$wsWeapon.Range("D28").Select()

# --- EnemyDB sheet: rebalance speed/health/damage/bulletSpeed ---
$wsEnemy = $wb.Worksheets.Item("EnemyDB")

$enemyRows = @(
    @(2,100,4,5,1),
    @(3,50,3,5,0.5),
    @(4,50,3,5,3),
    @(5,100,5,5,4)
)

foreach ($row in $enemyRows) {
    $r = $row[0]
    $wsEnemy.Cells.Item($r, 1).Value = $row[1]
    $wsEnemy.Cells.Item($r, 2).Value = $row[2]
    $wsEnemy.Cells.Item($r, 3).Value = $row[3]
    $wsEnemy.Cells.Item($r, 4).Value = $row[4]
}

$wsEnemy.Activate()
$wsEnemy.Range("E14").Select()
